$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 33   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/26/2026  Through  2/1/2026"

# --- Cells that become the special "0" / "***.*" text markers (copy format+value from an untouched template cell so style index s=13 is preserved) ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 200
$ws.Range("L15").Value = 20
$ws.Range("N15").Value = 100
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = -35.294117647058
$ws.Range("I16").Value = 24
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = -40
$ws.Range("L16").Value = -45.454545454545
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = -78.947368421052
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = 17.142857142857
$ws.Range("I17").Value = 49
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = 19.512195121951
$ws.Range("L17").Value = -12.5
$ws.Range("M17").Value = 113.04347826087
$ws.Range("N17").Value = -40.963855421686
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 22
$ws.Range("J18").Value = 23
$ws.Range("K18").Value = -4.347826086956
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 144.444444444444
$ws.Range("N18").Value = -70.666666666666
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -71.428571428571
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -45
$ws.Range("I19").Value = 27
$ws.Range("J19").Value = 46
$ws.Range("K19").Value = -41.304347826087
$ws.Range("L19").Value = -54.237288135593
$ws.Range("M19").Value = 58.823529411764
$ws.Range("N19").Value = -20.588235294117
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -52.941176470588
$ws.Range("I20").Value = 10
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -41.176470588235
$ws.Range("L20").Value = -56.521739130434
$ws.Range("M20").Value = -28.571428571428
$ws.Range("N20").Value = -79.166666666666
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -23.076923076923
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = -18.791946308724
$ws.Range("I21").Value = 139
$ws.Range("J21").Value = 169
$ws.Range("K21").Value = -17.751479289940
$ws.Range("L21").Value = -33.492822966507
$ws.Range("M21").Value = 63.529411764705
$ws.Range("N21").Value = -61.281337047353
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 17
$ws.Range("J23").Value = 16
$ws.Range("K23").Value = 6.25
$ws.Range("L23").Value = -56.410256410256
$ws.Range("M23").Value = -29.166666666666
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 20
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 88
$ws.Range("J24").Value = 93
$ws.Range("K24").Value = -5.376344086021
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 22.222222222222
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -66.666666666666
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 25
$ws.Range("I25").Value = 16
$ws.Range("J25").Value = 13
$ws.Range("K25").Value = 23.076923076923
$ws.Range("L25").Value = -11.111111111111
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -47.368421052631
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 65
$ws.Range("H26").Value = -33.846153846153
$ws.Range("I26").Value = 56
$ws.Range("J26").Value = 71
$ws.Range("K26").Value = -21.126760563380
$ws.Range("L26").Value = -18.840579710144
$ws.Range("M26").Value = -21.126760563380
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 100
$ws.Range("L27").Value = -25
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = 50
$ws.Range("L28").Value = 50
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = -50
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 100
$ws.Range("N30").Value = -50
